$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) store numeric-looking values ("319.92",
# "3.74%", ...) as plain text, not Number/Percentage cells. For each cell we
# touch in those columns, force a text number format first so Excel keeps the
# exact literal string instead of reinterpreting it as a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "319.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.74%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.54%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.247"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.36%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07743"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.59%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.693"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.28%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9440"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.00%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.26%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1241"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.72%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1854"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.59%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09236"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.83%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04324"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.04%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.45%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001294"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.56%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006003"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.83%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.343"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.11%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.339"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.39%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3333"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.56%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.774"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "11.20%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1354"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.85%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2825"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.34%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04036"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.21%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001269"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.34%"
$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004121"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.10%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.14%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02548"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "5.25%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007769"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.08%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1318"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.37%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007359"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.14%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.82%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008247"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "11.82%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.34%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006720"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.52%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.13%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2019"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "87.52%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004204"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "40.01%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.13%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.13%"
